$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$para1 = "Models for each of the 10 countries were retrained and serialized into JSON, saved into the models folder, and uploaded into Watson Studio."
$para2 = "Further code was developed to test out loading the JSON models, and using these to predict a 30 day forecast output for August 2019, and monthly and weekly trends per country" + [char]8217 + "s model."
$para3 = "The 30 day forecasts and trends per country can be found in a .docx file in the initial_model_predictions_trends folder."

$tr.Text = $para1 + "`r" + $para2 + "`r" + $para3

$italicWord1 = "models"
$offset1 = $para1.IndexOf("the " + $italicWord1 + " folder") + 4
$tr.Characters($offset1 + 1, $italicWord1.Length).Font.Italic = $true

$italicWord2 = "initial_model_predictions_trends"
$offset2 = $para1.Length + 1 + $para2.Length + 1 + $para3.IndexOf($italicWord2)
$tr.Characters($offset2 + 1, $italicWord2.Length).Font.Italic = $true
